# Applies the "Adicionando os textos as páginas" edit to Backlog.docx:
#  1. Ends the "Artes marciais..." paragraph with a period instead of a
#     trailing space ("...cunho esportivo " -> "...cunho esportivo.").
#  2. Swaps the god's name order: "Marte/Ares" -> "Ares/Marte".
#  3. Fixes the missing accent: "pratica" -> "prática".
#  4. Removes the "- Localização de academias próximas" bullet paragraph
#     from the "Desejável" backlog section.

$d = $word.ActiveDocument

# 1) "...cunho esportivo " -> "...cunho esportivo."
$d.Content.Find.Execute(
    "cunho esportivo ", $true, $false, $false, $false, $false,
    $true, 1, $false, "cunho esportivo.", 2) | Out-Null

# 2) "Marte/Ares ensinou" -> "Ares/Marte ensinou"
$d.Content.Find.Execute(
    "Marte/Ares ensinou", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ares/Marte ensinou", 2) | Out-Null

# 3) "para a pratica das artes marciais" -> "para a prática das artes marciais"
$d.Content.Find.Execute(
    "para a pratica das artes marciais", $true, $false, $false, $false, $false,
    $true, 1, $false, "para a prática das artes marciais", 2) | Out-Null

# 4) Remove the whole "- Localização de academias próximas" paragraph.
ForEach ($p in $d.Paragraphs) {
    If ($p.Range.Text.TrimEnd("`r`a") -eq "- Localização de academias próximas") {
        $p.Range.Delete()
        Break
    }
}
